$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pair table for Lrfn3-Lrfn3, updated per Dr Hou advice: added "ECs" cluster
# Full cartesian product of clusters {ECs, FAPs, sCs} x {ECs, FAPs, sCs} (rows 2-10)
$rows = @(
    @("ECs", "Lrfn3", "Lrfn3", "ECs", 2, 0.6666666666666666, 0.377313, 1.131939, 0.08069727054576252, 0.08069727054576252, 2, 0.6666666666666666, 0.377313, 1.131939, 0.08069727054576252, 0.08069727054576252, 0.142365099969, 1.281285899721, 0.006512049473535991, 0.006512049473535991),
    @("ECs", "Lrfn3", "Lrfn3", "FAPs", 2, 0.6666666666666666, 0.377313, 1.131939, 0.08069727054576252, 0.08069727054576252, 3, 1, 3.314087, 9.942261, 0.708795549719184, 0.7087955497191839, 1.250448108231, 11.254032974079, 0.05719786623732146, 0.05719786623732145),
    @("ECs", "Lrfn3", "Lrfn3", "sCs", 2, 0.6666666666666666, 0.377313, 1.131939, 0.08069727054576252, 0.08069727054576252, 3, 1, 0.98426, 2.95278, 0.2105071797350534, 0.2105071797350534, 0.37137409338, 3.34236684042, 0.01698735483490506, 0.01698735483490506),
    @("FAPs", "Lrfn3", "Lrfn3", "ECs", 3, 1, 3.314087, 9.942261, 0.708795549719184, 0.7087955497191839, 2, 0.6666666666666666, 0.377313, 1.131939, 0.08069727054576252, 0.08069727054576252, 1.250448108231, 11.254032974079, 0.05719786623732146, 0.05719786623732145),
    @("FAPs", "Lrfn3", "Lrfn3", "FAPs", 3, 1, 3.314087, 9.942261, 0.708795549719184, 0.7087955497191839, 3, 1, 3.314087, 9.942261, 0.708795549719184, 0.7087955497191839, 10.983172643569, 98.84855379212101, 0.5023911313017203, 0.5023911313017201),
    @("FAPs", "Lrfn3", "Lrfn3", "sCs", 3, 1, 3.314087, 9.942261, 0.708795549719184, 0.7087955497191839, 3, 1, 0.98426, 2.95278, 0.2105071797350534, 0.2105071797350534, 3.26192327062, 29.35730943558, 0.1492065521801423, 0.1492065521801423),
    @("sCs", "Lrfn3", "Lrfn3", "ECs", 3, 1, 0.98426, 2.95278, 0.2105071797350534, 0.2105071797350534, 2, 0.6666666666666666, 0.377313, 1.131939, 0.08069727054576252, 0.08069727054576252, 0.37137409338, 3.34236684042, 0.01698735483490506, 0.01698735483490506),
    @("sCs", "Lrfn3", "Lrfn3", "FAPs", 3, 1, 0.98426, 2.95278, 0.2105071797350534, 0.2105071797350534, 3, 1, 3.314087, 9.942261, 0.708795549719184, 0.7087955497191839, 3.26192327062, 29.35730943558, 0.1492065521801423, 0.1492065521801423),
    @("sCs", "Lrfn3", "Lrfn3", "sCs", 3, 1, 0.98426, 2.95278, 0.2105071797350534, 0.2105071797350534, 3, 1, 0.98426, 2.95278, 0.2105071797350534, 0.2105071797350534, 0.9687677476000001, 8.718909728400002, 0.04431327272000608, 0.04431327272000608)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
